$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "30.303.96"
Set-TextValue $ws.Cells.Item(2, 5) "  -0.39%  "

Set-TextValue $ws.Cells.Item(3, 4) "2.064.60"
Set-TextValue $ws.Cells.Item(3, 5) "  +3.18%  "

Set-TextValue $ws.Cells.Item(4, 5) "  +0.20%  "

Set-TextValue $ws.Cells.Item(5, 4) "326.14"
Set-TextValue $ws.Cells.Item(5, 5) "  +0.56%  "

Set-TextValue $ws.Cells.Item(6, 5) "  +0.19%  "

Set-TextValue $ws.Cells.Item(7, 4) "0.5167"
Set-TextValue $ws.Cells.Item(7, 5) "  +1.25%  "

Set-TextValue $ws.Cells.Item(8, 4) "0.4308"
Set-TextValue $ws.Cells.Item(8, 5) "  +4.22%  "

Set-TextValue $ws.Cells.Item(9, 4) "0.08675"
Set-TextValue $ws.Cells.Item(9, 5) "  -0.58%  "

Set-TextValue $ws.Cells.Item(10, 4) "45.64"
Set-TextValue $ws.Cells.Item(10, 5) "  +5.80%  "

Set-TextValue $ws.Cells.Item(11, 4) "1.146"
Set-TextValue $ws.Cells.Item(11, 5) "  +1.18%  "

Set-TextValue $ws.Cells.Item(12, 4) "23.97"
Set-TextValue $ws.Cells.Item(12, 5) "  -2.54%  "

Set-TextValue $ws.Cells.Item(13, 4) "2.066.66"
Set-TextValue $ws.Cells.Item(13, 5) "  +3.40%  "

Set-TextValue $ws.Cells.Item(14, 4) "6.583"
Set-TextValue $ws.Cells.Item(14, 5) "  +0.39%  "

Set-TextValue $ws.Cells.Item(15, 4) "7.601"
Set-TextValue $ws.Cells.Item(15, 5) "  +1.79%  "

Set-TextValue $ws.Cells.Item(16, 4) "1.005"
Set-TextValue $ws.Cells.Item(16, 5) "  +0.22%  "

Set-TextValue $ws.Cells.Item(17, 4) "94.41"
Set-TextValue $ws.Cells.Item(17, 5) "  +0.17%  "

Set-TextValue $ws.Cells.Item(18, 4) "0.00001112"
Set-TextValue $ws.Cells.Item(18, 5) "  -0.14%  "

Set-TextValue $ws.Cells.Item(19, 4) "0.06608"
Set-TextValue $ws.Cells.Item(19, 5) "  +1.61%  "

Set-TextValue $ws.Cells.Item(20, 4) "18.63"
Set-TextValue $ws.Cells.Item(20, 5) "  -1.32%  "

Set-TextValue $ws.Cells.Item(21, 5) "  +0.26%  "

Set-TextValue $ws.Cells.Item(22, 4) "6.166"
Set-TextValue $ws.Cells.Item(22, 5) "  -0.31%  "

Set-TextValue $ws.Cells.Item(23, 4) "30.346.15"
Set-TextValue $ws.Cells.Item(23, 5) "  -0.40%  "

Set-TextValue $ws.Cells.Item(24, 4) "12.16"
Set-TextValue $ws.Cells.Item(24, 5) "  +2.29%  "

Set-TextValue $ws.Cells.Item(25, 4) "2.273"
Set-TextValue $ws.Cells.Item(25, 5) "  +2.05%  "

Set-TextValue $ws.Cells.Item(26, 4) "2.313.97"
Set-TextValue $ws.Cells.Item(26, 5) "  +3.74%  "

Set-TextValue $ws.Cells.Item(27, 4) "22.02"
Set-TextValue $ws.Cells.Item(27, 5) "  -1.31%  "

Set-TextValue $ws.Cells.Item(28, 4) "160.27"
Set-TextValue $ws.Cells.Item(28, 5) "  -1.79%  "

Set-TextValue $ws.Cells.Item(29, 4) "2.483"
Set-TextValue $ws.Cells.Item(29, 5) "  +3.70%  "

Set-TextValue $ws.Cells.Item(30, 4) "129.94"
Set-TextValue $ws.Cells.Item(30, 5) "  -1.14%  "

Set-TextValue $ws.Cells.Item(31, 5) "  +2.64%  "

Set-TextValue $ws.Cells.Item(32, 4) "0.1060"
Set-TextValue $ws.Cells.Item(32, 5) "  +0.86%  "

Set-TextValue $ws.Cells.Item(33, 4) "6.024"
Set-TextValue $ws.Cells.Item(33, 5) "  -0.52%  "

Set-TextValue $ws.Cells.Item(34, 4) "3.835"
Set-TextValue $ws.Cells.Item(34, 5) "  -0.13%  "

Set-TextValue $ws.Cells.Item(35, 4) "1.485"
Set-TextValue $ws.Cells.Item(35, 5) "  +11.01%  "

Set-TextValue $ws.Cells.Item(36, 4) "0.02539"
Set-TextValue $ws.Cells.Item(36, 5) "  +0.92%  "

Set-TextValue $ws.Cells.Item(37, 4) "9.529"
Set-TextValue $ws.Cells.Item(37, 5) "  +5.64%  "

Set-TextValue $ws.Cells.Item(38, 4) "5.398"
Set-TextValue $ws.Cells.Item(38, 5) "  -0.58%  "

Set-TextValue $ws.Cells.Item(39, 4) "0.06552"
Set-TextValue $ws.Cells.Item(39, 5) "  -0.63%  "

Set-TextValue $ws.Cells.Item(40, 4) "12.36"
Set-TextValue $ws.Cells.Item(40, 5) "  -0.73%  "

Set-TextValue $ws.Cells.Item(41, 4) "0.2214"
Set-TextValue $ws.Cells.Item(41, 5) "  +0.86%  "

Set-TextValue $ws.Cells.Item(42, 4) "0.6592"
Set-TextValue $ws.Cells.Item(42, 5) "  -0.48%  "

Set-TextValue $ws.Cells.Item(43, 4) "1.229"
Set-TextValue $ws.Cells.Item(43, 5) "  -0.28%  "

Set-TextValue $ws.Cells.Item(44, 4) "1.003"
Set-TextValue $ws.Cells.Item(44, 5) "  +0.10%  "

Set-TextValue $ws.Cells.Item(45, 4) "13.91"
Set-TextValue $ws.Cells.Item(45, 5) "  +2.30%  "

Set-TextValue $ws.Cells.Item(46, 4) "0.6232"
Set-TextValue $ws.Cells.Item(46, 5) "  +1.21%  "

Set-TextValue $ws.Cells.Item(47, 4) "2.175"
Set-TextValue $ws.Cells.Item(47, 5) "  -0.90%  "

Set-TextValue $ws.Cells.Item(48, 4) "3.598"
Set-TextValue $ws.Cells.Item(48, 5) "  -1.85%  "

Set-TextValue $ws.Cells.Item(49, 4) "1.230"
Set-TextValue $ws.Cells.Item(49, 5) "  -2.86%  "

# Row 50: Aave -> WEMIXTOKEN
Set-TextValue $ws.Cells.Item(50, 2) "WEMIXTOKEN"
Set-TextValue $ws.Cells.Item(50, 3) "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Cells.Item(50, 4) "1.170"
Set-TextValue $ws.Cells.Item(50, 5) "  +5.62%  "

# Row 51: WEMIXTOKEN -> Aave
Set-TextValue $ws.Cells.Item(51, 2) "Aave"
Set-TextValue $ws.Cells.Item(51, 3) "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Cells.Item(51, 4) "80.73"
Set-TextValue $ws.Cells.Item(51, 5) "  +0.61%  "
